# Insert a new "cut" (ValueOrderDimension) into the Lookups sheet's cuts table,
# between the existing Gender_LIB cut (columns M:N) and the question_code cut
# (previously at O:P). This pushes question_code to Q:R, result_type to S:T,
# and zero_string to U1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# Insert two new blank columns at O:P; everything from O onward shifts right by 2.
$ws.Columns("O:P").Insert()

# Populate the new cut's header and value/count pairs.
$ws.Range("O1").Value = "ValueOrderDimension"
$ws.Range("O2").Value = "B"
$ws.Range("P2").Value = 2
$ws.Range("O3").Value = "A"
$ws.Range("P3").Value = 1

# The defined names that reference the tail of the cuts table need to be
# repointed to their new (shifted) locations.
$wb.Names.Item("cuts_head").RefersTo = "=Lookups!`$F`$1:`$T`$1"
$wb.Names.Item("zero_string").RefersTo = "=Lookups!`$U`$1"
